$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45041
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 24000
$ws.Range("L2").Value = 24000
$ws.Range("M2").Value = 24000
$ws.Range("N2").Value = '$/caja 18 kilos empedrada'
$ws.Range("P2").Value = 1333
$ws.Range("Q2").Value = 18

$ws.Range("D3").Value = 45037
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 24000
$ws.Range("L3").Value = 24000
$ws.Range("M3").Value = 24000
$ws.Range("N3").Value = '$/caja 15 kilos empedrada'
$ws.Range("P3").Value = 1600
$ws.Range("Q3").Value = 15

$ws.Range("D4").Value = 45034
$ws.Range("I4").Value = 'Primera'
$ws.Range("J4").Value = 50
$ws.Range("K4").Value = 24000
$ws.Range("L4").Value = 24000
$ws.Range("M4").Value = 24000
$ws.Range("N4").Value = '$/caja 18 kilos granel'
$ws.Range("P4").Value = 1333
$ws.Range("Q4").Value = 18

$ws.Range("D5").Value = 45033
$ws.Range("I5").Value = 'Primera'
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 24000
$ws.Range("L5").Value = 24000
$ws.Range("M5").Value = 24000
$ws.Range("N5").Value = '$/caja 18 kilos granel'
$ws.Range("P5").Value = 1333
$ws.Range("Q5").Value = 18

$ws.Range("D6").Value = 44285
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 20
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 25000
$ws.Range("N6").Value = '$/caja 18 kilos empedrada'
$ws.Range("P6").Value = 1389
$ws.Range("Q6").Value = 18

$ws.Range("D7").Value = 45043
$ws.Range("I7").Value = 'Primera'
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 24000
$ws.Range("L7").Value = 24000
$ws.Range("M7").Value = 24000
$ws.Range("N7").Value = '$/caja 18 kilos empedrada'
$ws.Range("P7").Value = 1333
$ws.Range("Q7").Value = 18

$ws.Range("D8").Value = 45042
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 24000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 24000
$ws.Range("N8").Value = '$/caja 18 kilos empedrada'
$ws.Range("P8").Value = 1333
$ws.Range("Q8").Value = 18

$ws.Range("D9").Value = 44313
$ws.Range("I9").Value = 'Primera'
$ws.Range("J9").Value = 20
$ws.Range("K9").Value = 15000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 15000
$ws.Range("N9").Value = '$/caja 15 kilos empedrada'
$ws.Range("P9").Value = 1000
$ws.Range("Q9").Value = 15

$ws.Range("D10").Value = 44313
$ws.Range("I10").Value = 'Primera'
$ws.Range("J10").Value = 20
$ws.Range("K10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("M10").Value = 30000
$ws.Range("N10").Value = '$/caja 20 kilos empedrada'
$ws.Range("P10").Value = 1500
$ws.Range("Q10").Value = 20

$ws.Range("D11").Value = 45015
$ws.Range("I11").Value = 'Primera'
$ws.Range("J11").Value = 50
$ws.Range("K11").Value = 24000
$ws.Range("L11").Value = 24000
$ws.Range("M11").Value = 24000
$ws.Range("N11").Value = '$/caja 18 kilos empedrada'
$ws.Range("P11").Value = 1333
$ws.Range("Q11").Value = 18

$ws.Range("D12").Value = 45044
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 40
$ws.Range("K12").Value = 24000
$ws.Range("L12").Value = 24000
$ws.Range("M12").Value = 24000
$ws.Range("N12").Value = '$/caja 18 kilos empedrada'
$ws.Range("P12").Value = 1333
$ws.Range("Q12").Value = 18

$ws.Range("D13").Value = 44293
$ws.Range("I13").Value = 'Primera'
$ws.Range("J13").Value = 10
$ws.Range("K13").Value = 25000
$ws.Range("L13").Value = 25000
$ws.Range("M13").Value = 25000
$ws.Range("N13").Value = '$/caja 15 kilos empedrada'
$ws.Range("P13").Value = 1667
$ws.Range("Q13").Value = 15

$ws.Range("D14").Value = 45014
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 8000
$ws.Range("L14").Value = 8000
$ws.Range("M14").Value = 8000
$ws.Range("N14").Value = '$/caja 18 kilos empedrada'
$ws.Range("P14").Value = 444
$ws.Range("Q14").Value = 18

$ws.Range("D15").Value = 44315
$ws.Range("I15").Value = 'Especial'
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 30000
$ws.Range("L15").Value = 30000
$ws.Range("M15").Value = 30000
$ws.Range("N15").Value = '$/caja 20 kilos empedrada'
$ws.Range("P15").Value = 1500
$ws.Range("Q15").Value = 20

$ws.Range("D16").Value = 44315
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = 15000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 15000
$ws.Range("N16").Value = '$/caja 15 kilos granel'
$ws.Range("P16").Value = 1000
$ws.Range("Q16").Value = 15

$ws.Range("D17").Value = 44280
$ws.Range("I17").Value = 'Primera'
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 25000
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = 25000
$ws.Range("N17").Value = '$/caja 18 kilos empedrada'
$ws.Range("P17").Value = 1389
$ws.Range("Q17").Value = 18

$ws.Range("D18").Value = 45040
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 80
$ws.Range("K18").Value = 18000
$ws.Range("L18").Value = 18000
$ws.Range("M18").Value = 18000
$ws.Range("N18").Value = '$/caja 18 kilos empedrada'
$ws.Range("P18").Value = 1000
$ws.Range("Q18").Value = 18
